$d = $word.ActiveDocument

function New-WordXml([string]$bodyFragment) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1) "Applicant_bkbtn" paragraph (right after Applicant_txtbx2, before the blank /
#    "*Employer_Login" paragraph): drop its <w:lastRenderedPageBreak/> -- the page
#    break marker moves down to the new "*Employer_ regis" paragraph instead.
$pBkBtn = $d.Paragraphs.Item(17)
$pBkBtn.Range.InsertXML((New-WordXml '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Applicant_bkbtn</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')) | Out-Null

# 2) Turn the "*Employer_Login" paragraph into the new "*Employer_ regis" heading
#    (with the page-break marker now here), then splice in the freshly written
#    Employer-registration UI stub paragraphs, and finally restore the original
#    "*Employer_Login" paragraph right after them (it still introduces the
#    pre-existing Employer login block further below).
$pEmployerLogin = $d.Paragraphs.Item(19)
$newBlock = ''
$newBlock += '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">*Employer_ </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>regis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>'
$newBlock += '<w:p><w:r><w:t>Employer_pnl1</w:t></w:r></w:p>'
$newBlock += '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Employer</w:t></w:r><w:r><w:t>_bk</w:t></w:r><w:r><w:t>btn</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p>'
$newBlock += '<w:p/>'
$newBlock += '<w:p/>'
$newBlock += '<w:p><w:r><w:t>*</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Employer_Login</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$pEmployerLogin.Range.InsertXML((New-WordXml $newBlock)) | Out-Null

# The block above replaced 1 paragraph with 6, so everything after it shifted
# down by 5 paragraphs.
$shift = 5

# 3) The page-break marker that used to sit on "Applicant_pcbx1" now renders on
#    the preceding "Applicant_lbl1" paragraph instead.
$pApplicantLbl1 = $d.Paragraphs.Item(32 + $shift)
$pApplicantLbl1.Range.InsertXML((New-WordXml '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Applicant_lbl1</w:t></w:r></w:p>')) | Out-Null

$pApplicantPcbx1 = $d.Paragraphs.Item(33 + $shift)
$pApplicantPcbx1.Range.InsertXML((New-WordXml '<w:p><w:r><w:t>Applicant_pcbx1</w:t></w:r></w:p>')) | Out-Null

# 4) The "_GoBack" bookmark that used to sit on the final "Employer..._pcbx1"
#    paragraph moved up to the new "Employer_bkbtn" paragraph (step 2), so drop
#    it from here.
$pEmployerPcbx1 = $d.Paragraphs.Item(44 + $shift)
$pEmployerPcbx1.Range.InsertXML((New-WordXml '<w:p><w:r><w:t>Employer</w:t></w:r><w:r><w:t>_pcbx1</w:t></w:r></w:p>')) | Out-Null
